# Fruta / hortaliza, semanal
# Adds three new weekly price records to the "Fruta, Feria Lagunitas de
# Puerto Montt - Frutilla" dataset. The existing rows 283:355 (Fecha
# serial 44817 .. 44306) shift down to 286:358, and three brand-new
# rows are inserted at 283:285 with Fecha = 44943 (2023-01-17).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows right before the current row 283; this pushes the
# old rows 283:355 down to 286:358, matching the rest of the diff
# (which is just every old row's data shifted down by three).
$ws.Range("A283:T285").EntireRow.Insert()

# ---- New row 283 ----
$ws.Range("A283").Value = 4
$ws.Range("B283").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C283").Value = "Los Lagos"
$ws.Range("D283").Value = 44943
$ws.Range("E283").Value = 10
$ws.Range("F283").Value = "Fruta"
$ws.Range("G283").Value = 100101
$ws.Range("H283").Value = "Berries"
$ws.Range("I283").Value = 100112025
$ws.Range("J283").Value = "Frutilla"
$ws.Range("K283").Value = "Sin especificar"
$ws.Range("L283").Value = "Primera"
$ws.Range("M283").Value = 200
$ws.Range("N283").Value = 10000
$ws.Range("O283").Value = 11000
$ws.Range("P283").Value = 10500
$ws.Range("Q283").Value = "`$/bandeja 7 kilos"
$ws.Range("R283").Value = "Provincia de Melipilla"
$ws.Range("S283").Value = 1500
$ws.Range("T283").Value = 7

# ---- New row 284 ----
$ws.Range("A284").Value = 4
$ws.Range("B284").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C284").Value = "Los Lagos"
$ws.Range("D284").Value = 44943
$ws.Range("E284").Value = 10
$ws.Range("F284").Value = "Fruta"
$ws.Range("G284").Value = 100101
$ws.Range("H284").Value = "Berries"
$ws.Range("I284").Value = 100112025
$ws.Range("J284").Value = "Frutilla"
$ws.Range("K284").Value = "Sin especificar"
$ws.Range("L284").Value = "Primera"
$ws.Range("M284").Value = 600
$ws.Range("N284").Value = 9000
$ws.Range("O284").Value = 10000
$ws.Range("P284").Value = 9500
$ws.Range("Q284").Value = "`$/caja 7 kilos"
$ws.Range("R284").Value = "Región de La Araucanía"
$ws.Range("S284").Value = 1357
$ws.Range("T284").Value = 7

# ---- New row 285 ----
$ws.Range("A285").Value = 4
$ws.Range("B285").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C285").Value = "Los Lagos"
$ws.Range("D285").Value = 44943
$ws.Range("E285").Value = 10
$ws.Range("F285").Value = "Fruta"
$ws.Range("G285").Value = 100101
$ws.Range("H285").Value = "Berries"
$ws.Range("I285").Value = 100112025
$ws.Range("J285").Value = "Frutilla"
$ws.Range("K285").Value = "Sin especificar"
$ws.Range("L285").Value = "Segunda"
$ws.Range("M285").Value = 300
$ws.Range("N285").Value = 8000
$ws.Range("O285").Value = 8000
$ws.Range("P285").Value = 8000
$ws.Range("Q285").Value = "`$/caja 7 kilos"
$ws.Range("R285").Value = "Región de La Araucanía"
$ws.Range("S285").Value = 1143
$ws.Range("T285").Value = 7
